$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I14").NumberFormat = '#,##0'
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = -50
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -66.666666666666
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 1
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("F15").Value = 1
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = 60
$ws.Range("M15").Value = -52.941176470588
$ws.Range("N15").Value = -52.941176470588
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 51
$ws.Range("K16").Value = 59.375
$ws.Range("L16").Value = 18.60465116279
$ws.Range("M16").Value = -40.697674418604
$ws.Range("N16").Value = -79.183673469387
$ws.Range("C17").Value = 1
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").NumberFormat = '@'
$ws.Range("E17").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 142.857142857143
$ws.Range("I17").Value = 109
$ws.Range("K17").Value = 3.809523809523
$ws.Range("L17").Value = 22.471910112359
$ws.Range("M17").Value = -6.03448275862
$ws.Range("N17").Value = -59.328358208955
$ws.Range("C18").NumberFormat = '@'
$ws.Range("C18").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("D18").Value = 1
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = 100
$ws.Range("J18").Value = 63
$ws.Range("K18").Value = 3.174603174603
$ws.Range("L18").Value = -13.333333333333
$ws.Range("M18").Value = -61.538461538461
$ws.Range("N18").Value = -94.273127753304
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 14.285714285714
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 284
$ws.Range("J19").Value = 265
$ws.Range("K19").Value = 7.169811320754
$ws.Range("L19").Value = 44.162436548223
$ws.Range("M19").Value = -17.681159420289
$ws.Range("N19").Value = -58.840579710144
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 9.090909090909
$ws.Range("I20").Value = 108
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 80
$ws.Range("L20").Value = 89.473684210526
$ws.Range("M20").Value = 9.090909090909
$ws.Range("N20").Value = -95.336787564766
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 60
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 49
$ws.Range("H21").Value = 34.69387755102
$ws.Range("I21").Value = 626
$ws.Range("J21").Value = 533
$ws.Range("K21").Value = 17.448405253283
$ws.Range("L21").Value = 34.047109207708
$ws.Range("M21").Value = -24.84993997599
$ws.Range("N21").Value = -86.60676080445
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("C23").Value = 1
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E23").Value = 0
$ws.Range("F23").NumberFormat = '#,##0'
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 16
$ws.Range("J23").Value = 34
$ws.Range("K23").Value = -52.941176470588
$ws.Range("L23").Value = -5.882352941176
$ws.Range("M23").Value = -15.78947368421
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 130
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 51
$ws.Range("H24").Value = 98.039215686274
$ws.Range("I24").Value = 711
$ws.Range("J24").Value = 404
$ws.Range("K24").Value = 75.990099009901
$ws.Range("L24").Value = 39.138943248532
$ws.Range("M24").Value = -50.314465408805
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 133.333333333333
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = -31.818181818181
$ws.Range("I25").Value = 280
$ws.Range("J25").Value = 207
$ws.Range("K25").Value = 35.265700483091
$ws.Range("L25").Value = 43.589743589743
$ws.Range("M25").Value = -38.461538461538
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("C26").Value = 1
$ws.Range("F26").NumberFormat = '#,##0'
$ws.Range("F26").Value = 1
$ws.Range("I26").Value = 17
$ws.Range("K26").Value = 30.76923076923
$ws.Range("L26").Value = 30.76923076923
$ws.Range("D27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = -3.703703703703
$ws.Range("L27").Value = 52.941176470588
# --- Header text updates (shared-string rich-text runs) ---

# A8: "Volume 29   Number  42" -> "...43"
$cellA8 = $ws.Range("A8")
$textA8 = $cellA8.Value2
$posA8 = $textA8.IndexOf("42")
$charsA8 = $cellA8.Characters($posA8 + 1, 2)
$charsA8.Text = "43"

# C9: "Report Covering the Week  10/17/2022  Through  10/23/2022"
#     -> "...10/24/2022  Through  10/30/2022"
$cellC9 = $ws.Range("C9")
$textC9 = $cellC9.Value2
$pos1C9 = $textC9.IndexOf("10/17/2022")
$chars1C9 = $cellC9.Characters($pos1C9 + 1, 10)
$chars1C9.Text = "10/24/2022"

$textC9b = $cellC9.Value2
$pos2C9 = $textC9b.IndexOf("10/23/2022")
$chars2C9 = $cellC9.Characters($pos2C9 + 1, 10)
$chars2C9.Text = "10/30/2022"
